$wb = $excel.ActiveWorkbook

# --- "studies" sheet (sheet1): add PMID column in H ---
$wsStudies = $wb.Worksheets.Item("studies")
$wsStudies.Range("H1").Value = "PMID"

# --- "counts" sheet (sheet3): add notes column in F ---
$wsCounts = $wb.Worksheets.Item("counts")
$wsCounts.Range("F1").Value = "notes"

# Update selections to match the recorded user activity after the edit
$wsStudies.Range("H2").Select()
$wsCounts.Range("F2").Select()

# "counts" sheet was the active tab
$wsCounts.Activate()
